$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.295.21"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "1.831.48"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  +0.79%  "

$ws.Range("D5").Value = "'314.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").Value = "'0.4740"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.3679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("E9").Value = "  +1.13%  "

$ws.Range("D10").Value = "'0.8849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").Value = "1.905.87"
$ws.Range("E12").Value = "  +1.93%  "

$ws.Range("D13").Value = "'0.07307"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("D14").Value = "'5.419"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").Value = "'93.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "

$ws.Range("D16").Value = "'6.551"
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "'0.000008795"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").Value = "27.728.87"
$ws.Range("E20").Value = "  +2.76%  "

$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "2.124.60"
$ws.Range("E24").Value = "  +3.94%  "

$ws.Range("D25").Value = "'1.904"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").Value = "'151.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("D28").Value = "'2.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").Value = "'5.230"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").Value = "'117.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "'0.08980"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7487"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.174"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").Value = "'4.538"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("D35").Value = "'2.946"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").Value = "'1.095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").Value = "'0.01954"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").Value = "'2.417"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "

$ws.Range("D41").Value = "'2.948"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("D42").Value = "'7.223"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "'0.5288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").Value = "'0.1656"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "'8.484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("D46").Value = "'0.4905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").Value = "'10.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "'105.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").Value = "'1.663"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
